# The deck's design/theme was changed from the custom "Integral" theme to
# the default Office "Office Theme" (i.e. the presentation's applied theme
# colours were reset to the stock Office palette, as the earlier "Office
# Theme" colours were already stashed in the deck's secondary theme slot).
#
# The ThemeColorScheme on the (single) SlideMaster's Theme is the editable,
# persisted surface for the applied theme's colour scheme, in the fixed
# COM order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

# Target palette: the stock PowerPoint "Office Theme" colour scheme.
# RGB values are encoded PowerPoint/OLE-style as 0x00BBGGRR (i.e. R + G*256 + B*65536).
$colors.Item(1).RGB  = 0            # dk1       000000
$colors.Item(2).RGB  = 16777215     # lt1       FFFFFF
$colors.Item(3).RGB  = 6968388      # dk2       44546A
$colors.Item(4).RGB  = 15132391     # lt2       E7E6E6
$colors.Item(5).RGB  = 13998939     # accent1   5B9BD5
$colors.Item(6).RGB  = 3243501      # accent2   ED7D31
$colors.Item(7).RGB  = 10855845     # accent3   A5A5A5
$colors.Item(8).RGB  = 49407        # accent4   FFC000
$colors.Item(9).RGB  = 12874308     # accent5   4472C4
$colors.Item(10).RGB = 4697456      # accent6   70AD47
$colors.Item(11).RGB = 12673797     # hlink     0563C1
$colors.Item(12).RGB = 7491477      # folHlink  954F72
